$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "59.700.51"
$ws.Range("E2").Value2 = "  -1.96%  "
$ws.Range("D3").Value2 = "2.589.15"
$ws.Range("E3").Value2 = "  -3.17%  "
$ws.Range("E4").Value2 = "  +0.11%  "
$ws.Range("D5").Value2 = "'561.14"
$ws.Range("E5").Value2 = "  -1.37%  "
$ws.Range("D6").Value2 = "'143.08"
$ws.Range("E6").Value2 = "  -2.89%  "
$ws.Range("D7").Value2 = "'0.999"
$ws.Range("E7").Value2 = "  +0.09%  "
$ws.Range("D8").Value2 = "'0.599"
$ws.Range("E8").Value2 = "  -1.41%  "
$ws.Range("D9").Value2 = "2.599.67"
$ws.Range("E9").Value2 = "  -2.70%  "
$ws.Range("E10").Value2 = "  -2.93%  "
$ws.Range("E11").Value2 = "  -0.80%  "
$ws.Range("D12").Value2 = "'0.162"
$ws.Range("E12").Value2 = "  +11.64%  "
$ws.Range("E13").Value2 = "  +4.01%  "
$ws.Range("D14").Value2 = "3.045.76"
$ws.Range("E14").Value2 = "  -2.23%  "
$ws.Range("B15").Value2 = "WrappedBTC"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value2 = "59.639.58"
$ws.Range("E15").Value2 = "  -1.88%  "
$ws.Range("B16").Value2 = "Avalanche"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value2 = "'23.29"
$ws.Range("E16").Value2 = "  +5.77%  "
$ws.Range("E17").Value2 = "  -0.40%  "
$ws.Range("D18").Value2 = "2.590.82"
$ws.Range("E18").Value2 = "  -2.50%  "
$ws.Range("D19").Value2 = "'4.59"
$ws.Range("E19").Value2 = "  +0.18%  "
$ws.Range("D20").Value2 = "'338.97"
$ws.Range("E20").Value2 = "  -1.50%  "
$ws.Range("D21").Value2 = "'10.42"
$ws.Range("E21").Value2 = "  -0.71%  "
$ws.Range("D22").Value2 = "'6.52"
$ws.Range("E22").Value2 = "  +2.21%  "
$ws.Range("E23").Value2 = "  +0.32%  "
$ws.Range("D24").Value2 = "'63.57"
$ws.Range("E24").Value2 = "  -5.09%  "
$ws.Range("D25").Value2 = "'0.474"
$ws.Range("E26").Value2 = "  +0.15%  "
$ws.Range("D27").Value2 = "'0.162"
$ws.Range("E27").Value2 = "  -2.02%  "
$ws.Range("D28").Value2 = "'7.48"
$ws.Range("E28").Value2 = "  +0.56%  "
$ws.Range("D29").Value2 = "0.0₃0781"
$ws.Range("E29").Value2 = "  -4.25%  "
$ws.Range("E30").Value2 = "  +0.08%  "
$ws.Range("D31").Value2 = "'6.21"
$ws.Range("E31").Value2 = "  -0.63%  "
$ws.Range("E32").Value2 = "  -2.18%  "
$ws.Range("D33").Value2 = "'157.93"
$ws.Range("E33").Value2 = "  +0.96%  "
$ws.Range("D34").Value2 = "'19.11"
$ws.Range("E34").Value2 = "  -0.82%  "
$ws.Range("E35").Value2 = "  -1.05%  "
$ws.Range("E36").Value2 = "  +1.00%  "
$ws.Range("E37").Value2 = "  -0.57%  "
$ws.Range("D38").Value2 = "'0.867"
$ws.Range("E38").Value2 = "  -4.74%  "
$ws.Range("D39").Value2 = "'37.36"
$ws.Range("E39").Value2 = "  -0.61%  "
$ws.Range("E40").Value2 = "  -2.01%  "
$ws.Range("D41").Value2 = "'295.45"
$ws.Range("E41").Value2 = "  -2.83%  "
$ws.Range("D42").Value2 = "'3.69"
$ws.Range("E42").Value2 = "  +0.41%  "
$ws.Range("D43").Value2 = "'137.73"
$ws.Range("E43").Value2 = "  +6.93%  "
$ws.Range("D44").Value2 = "'1.00"
$ws.Range("E44").Value2 = "  +0.08%  "
$ws.Range("D45").Value2 = "'0.0978"
$ws.Range("E45").Value2 = "  -0.73%  "
$ws.Range("D46").Value2 = "'0.596"
$ws.Range("E46").Value2 = "  -2.02%  "
$ws.Range("D47").Value2 = "'10.65"
$ws.Range("E47").Value2 = "  -0.63%  "
$ws.Range("D48").Value2 = "'0.0533"
$ws.Range("E48").Value2 = "  -2.79%  "
$ws.Range("E49").Value2 = "  -0.60%  "
$ws.Range("B50").Value2 = "InjectiveProtocol"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value2 = "'18.78"
$ws.Range("E50").Value2 = "  -1.03%  "
$ws.Range("B51").Value2 = "Maker"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value2 = "1.971.05"
$ws.Range("E51").Value2 = "  +0.15%  "
